$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (pushing existing rows 4-32 down to 5-33),
# then populate it with a new week's data record (same values as the
# former last row, but with a new, more recent date).
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44817
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112044
$ws.Range("G4").Value = "Perejil"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 900
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = 950
$ws.Range("N4").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O4").Value = "Región de Arica y Parinacota"
$ws.Range("P4").Value = 475
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = "Hortaliza"
